# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and E (Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell without letting Excel re-interpret numeric-looking
# strings (e.g. "1.001") as a floating point number, and without leaving behind any
# number-format / style change on the cell.
function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "28.524.94"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.823.04"
$ws.Range("E3").Value = "  -0.22%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue "D5" "315.16"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  +0.03%  "
Set-TextValue "D7" "0.5106"
$ws.Range("E7").Value = "  -6.07%  "
Set-TextValue "D8" "0.3957"
$ws.Range("E8").Value = "  -1.97%  "
Set-TextValue "D9" "0.08173"
$ws.Range("E9").Value = "  +6.43%  "
Set-TextValue "D10" "1.111"
$ws.Range("E10").Value = "  -0.82%  "
Set-TextValue "D11" "41.67"
$ws.Range("E11").Value = "  -0.49%  "
Set-TextValue "D12" "21.14"
$ws.Range("E12").Value = "  +0.23%  "
Set-TextValue "D13" "6.321"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +0.07%  "
Set-TextValue "D15" "7.525"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "1.825.06"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  +3.68%  "
Set-TextValue "D18" "92.67"
$ws.Range("E18").Value = "  +3.00%  "
Set-TextValue "D19" "0.06658"
$ws.Range("E19").Value = "  +0.76%  "
Set-TextValue "D20" "17.82"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.04%  "
Set-TextValue "D22" "6.095"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "28.556.93"
$ws.Range("E23").Value = "  +0.08%  "
Set-TextValue "D24" "11.42"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  -0.56%  "
Set-TextValue "D26" "21.48"
$ws.Range("E26").Value = "  +3.32%  "
Set-TextValue "D27" "156.60"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "2.033.34"
Set-TextValue "D29" "2.399"
$ws.Range("E29").Value = "  -2.13%  "
Set-TextValue "D30" "126.45"
$ws.Range("E30").Value = "  +1.94%  "
Set-TextValue "D31" "1.110"
$ws.Range("E31").Value = "  -1.33%  "
Set-TextValue "D32" "0.1091"
$ws.Range("E32").Value = "  -1.40%  "
Set-TextValue "D33" "5.750"
$ws.Range("E33").Value = "  +1.23%  "
Set-TextValue "D34" "3.662"
$ws.Range("E34").Value = "  +0.48%  "
Set-TextValue "D35" "0.07045"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("E36").Value = "  -0.55%  "
Set-TextValue "D37" "5.274"
$ws.Range("E37").Value = "  +1.15%  "
Set-TextValue "D38" "0.02352"
$ws.Range("E38").Value = "  -0.12%  "
Set-TextValue "D39" "8.836"
$ws.Range("E39").Value = "  -0.73%  "
Set-TextValue "D40" "0.6323"
$ws.Range("E40").Value = "  +0.37%  "
Set-TextValue "D41" "11.30"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -0.55%  "
Set-TextValue "D43" "1.000"
$ws.Range("E43").Value = "  +0.04%  "
Set-TextValue "D44" "1.398"
$ws.Range("E44").Value = "  -0.24%  "
Set-TextValue "D45" "13.53"
$ws.Range("E45").Value = "  +0.77%  "
Set-TextValue "D46" "0.5942"
$ws.Range("E46").Value = "  +1.07%  "
Set-TextValue "D47" "3.733"
$ws.Range("E47").Value = "  +0.72%  "
Set-TextValue "D48" "125.42"
$ws.Range("E48").Value = "  +0.03%  "
Set-TextValue "D49" "1.996"
$ws.Range("E49").Value = "  -0.37%  "
Set-TextValue "D50" "1.190"
$ws.Range("E50").Value = "  -0.73%  "
Set-TextValue "D51" "0.06916"
$ws.Range("E51").Value = "  +0.13%  "
